$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.000.30"
$ws.Range("E2").Value = "  -0.75%  "
$ws.Range("D3").Value = "2.551.17"
$ws.Range("E3").Value = "  -0.12%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.64"
$ws.Range("E5").Value = "  +1.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.16"
$ws.Range("E6").Value = "  -1.96%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  -0.67%  "
$ws.Range("E9").Value = "  -0.93%  "
$ws.Range("E10").Value = "  -4.38%  "
$ws.Range("E11").Value = "  -0.35%  "
$ws.Range("E12").Value = "  -1.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.16"
$ws.Range("E13").Value = "  -3.17%  "
$ws.Range("D14").Value = "3.006.05"
$ws.Range("E14").Value = "  -0.22%  "
$ws.Range("D15").Value = "62.910.45"
$ws.Range("E15").Value = "  -0.83%  "
$ws.Range("E16").Value = "  -0.65%  "
$ws.Range("D17").Value = "2.556.66"
$ws.Range("E17").Value = "  -0.13%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.33"
$ws.Range("E18").Value = "  -2.27%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "335.22"
$ws.Range("E19").Value = "  -1.86%  "
$ws.Range("E20").Value = "  -0.50%  "
$ws.Range("E21").Value = "  -2.32%  "
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.52"
$ws.Range("E23").Value = "  -0.96%  "
$ws.Range("E24").Value = "  -0.37%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.60"
$ws.Range("E25").Value = "  +1.71%  "
$ws.Range("B26").Value = "SuiNetwork"
$ws.Range("C26").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.48"
$ws.Range("E26").Value = "  +0.56%  "
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.995"
$ws.Range("E27").Value = "  -0.58%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.37"
$ws.Range("E28").Value = "  -0.27%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.31"
$ws.Range("E29").Value = "  +1.35%  "
$ws.Range("E30").Value = "  +1.70%  "
$ws.Range("D31").Value = "0.0₃0812"
$ws.Range("E31").Value = "  -3.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "177.34"
$ws.Range("E32").Value = "  +0.05%  "
$ws.Range("E33").Value = "  -1.56%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "410.71"
$ws.Range("E34").Value = "  -0.80%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "19.13"
$ws.Range("E35").Value = "  +0.27%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.399"
$ws.Range("E36").Value = "  -1.82%  "
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("E38").Value = "  -2.60%  "
$ws.Range("E39").Value = "  -0.70%  "
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("E41").Value = "  -1.22%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "151.25"
$ws.Range("E42").Value = "  -2.78%  "
$ws.Range("E43").Value = "  -1.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.79"
$ws.Range("E44").Value = "  -1.25%  "
$ws.Range("E45").Value = "  +0.89%  "
$ws.Range("E46").Value = "  -1.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0966"
$ws.Range("E47").Value = "  +0.19%  "
$ws.Range("E48").Value = "  +1.96%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.24"
$ws.Range("E49").Value = "  -2.66%  "
$ws.Range("E50").Value = "  -8.73%  "
$ws.Range("E51").Value = "  -0.06%  "
